# Replace the schedule (columns B:G, rows 2-29) with a Latin-square schedule.
# Rows 26-29 are brand-new rows appended below the previous last row (25);
# column A (Proband-ID) is intentionally left blank for those new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    (1,4,6,5,4,3),
    (5,3,4,3,1,6),
    (6,3,2,4,5,1),
    (2,5,1,6,4,3),
    (3,6,5,1,2,4),
    (4,1,3,2,6,5),
    (1,4,6,5,4,3),
    (5,3,4,3,1,6),
    (6,3,2,4,5,1),
    (2,5,1,6,4,3),
    (3,6,5,1,2,4),
    (1,4,6,5,4,3),
    (5,3,4,3,1,6),
    (6,3,2,4,5,1),
    (2,5,1,6,4,3),
    (3,6,5,1,2,4),
    (1,4,6,5,4,3),
    (5,3,4,3,1,6),
    (6,3,2,4,5,1),
    (2,5,1,6,4,3),
    (3,6,5,1,2,4),
    (4,1,3,2,6,5),
    (1,4,6,5,4,3),
    (5,3,4,3,1,6),
    (6,3,2,4,5,1),
    (2,5,1,6,4,3),
    (3,6,5,1,2,4),
    (4,1,3,2,6,5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Update the view: scroll so row 7 is the top-left row, and select B28
# (matches the state the workbook was left in after the edit).
$ws.Range("B28").Select()
$excel.ActiveWindow.ScrollRow = 7
